# Apply crypto price/volume updates matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.543.51"
$ws.Range("E2").Value = "  +1.91%  "

$ws.Range("D3").Value = "2.601.10"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("D9").Value = "2.624.93"
$ws.Range("E9").Value = "  +1.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "

$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.157"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.371"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.35%  "

$ws.Range("D14").Value = "3.072.40"
$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.78%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "60.535.20"
$ws.Range("E16").Value = "  +1.92%  "

$ws.Range("E17").Value = "  +3.48%  "

$ws.Range("D18").Value = "2.622.92"
$ws.Range("E18").Value = "  +1.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.05%  "

$ws.Range("E22").Value = "  +7.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  +8.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("E28").Value = "  +6.56%  "

$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  +3.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.91%  "

$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.03%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.997"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "164.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.38%  "

$ws.Range("E34").Value = "  +2.43%  "

$ws.Range("E35").Value = "  +12.20%  "

$ws.Range("E36").Value = "  +3.92%  "

$ws.Range("E37").Value = "  +5.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.44%  "

$ws.Range("E40").Value = "  +6.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "310.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.17%  "

$ws.Range("E42").Value = "  -1.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0989"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.11%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.80%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.78%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.607"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0552"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.24%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.61%  "

$ws.Range("E51").Value = "  +3.59%  "
